$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.355.86"
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("D3").Value = "3.171.91"
$ws.Range("E3").Value = "  -3.68%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "612.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.387"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.684"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").Value = "3.180.09"
$ws.Range("E10").Value = "  -3.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.567"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.176"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.12%  "

$ws.Range("D14").Value = "3.862.39"

$ws.Range("D15").Value = "90.215.50"
$ws.Range("E15").Value = "  +1.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.94%  "

$ws.Range("D18").Value = "3.185.51"
$ws.Range("E18").Value = "  -3.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "431.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000188"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +38.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.96%  "

$ws.Range("D27").Value = "3.432.38"
$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "74.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.167"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.94%  "

$ws.Range("E31").Value = "  +0.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +31.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "532.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.85%  "

$ws.Range("E37").Value = "  -11.61%  "

$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.95%  "

$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.125"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.13%  "

$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.373"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "146.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "44.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "172.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.94%  "

$ws.Range("E48").Value = "  -3.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.610"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.88%  "
